# Generate Report for Handback
# - Status column flips from "Ready for handoff" to "Handed back: in sync with en-US"
#   for every localization row (zh-cn + de-de sheets).
# - Each row now has a populated "Latest Target File" (hyperlinked to the source
#   .md on GitHub, same target as column A), "Latest Handback File" (the handback
#   .xlf name) and a real "Latest Handback DateTime" (previously the epoch
#   placeholder "0001-01-01 00:00:00").

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

$urlMd32 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d3289b63a9cf6e8d5136b5e60614da12468ffb43/e2e/32e4e091-9df4-414b-bc4c-8236a1306f0c.md"
$urlMdD2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d3289b63a9cf6e8d5136b5e60614da12468ffb43/e2e/d2abd64f-9c0f-4af3-acf8-f4f5f601c048.md"

# ---------------- zh-cn ----------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("C2").Value = $statusText
$ws.Range("C3").Value = $statusText

$ws.Range("J2").Value = "32e4e091-9df4-414b-bc4c-8236a1306f0c.811ab3495dff9424ce59152fc652fd16cd023d27.zh-cn.xlf"
$ws.Range("K2").Value = "2016-08-31 05:40:38"

$ws.Range("J3").Value = "d2abd64f-9c0f-4af3-acf8-f4f5f601c048.05a88b359f114698fe9368b5815e5d2a3493305a.zh-cn.xlf"
$ws.Range("K3").Value = "2016-08-31 05:40:38"

# Recreate all four hyperlinks (A2, I2, A3, I3) in reading order so the
# relationship ids line up the way Excel would renumber them.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), $urlMd32, "", "", "32e4e091-9df4-414b-bc4c-8236a1306f0c.md")
$ws.Hyperlinks.Add($ws.Range("I2"), $urlMd32, "", "", "32e4e091-9df4-414b-bc4c-8236a1306f0c.md")
$ws.Hyperlinks.Add($ws.Range("A3"), $urlMdD2, "", "", "d2abd64f-9c0f-4af3-acf8-f4f5f601c048.md")
$ws.Hyperlinks.Add($ws.Range("I3"), $urlMdD2, "", "", "d2abd64f-9c0f-4af3-acf8-f4f5f601c048.md")

$ws.Range("I2").Font.Underline = $true
$ws.Range("I2").Font.Color = 15570276
$ws.Range("I3").Font.Underline = $true
$ws.Range("I3").Font.Color = 15570276

$ws.Columns.Item(3).ColumnWidth = 29.9777050018311
$ws.Columns.Item(9).ColumnWidth = 40
$ws.Columns.Item(10).ColumnWidth = 40

# ---------------- de-de ----------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("C2").Value = $statusText
$ws.Range("C3").Value = $statusText

$ws.Range("J2").Value = "32e4e091-9df4-414b-bc4c-8236a1306f0c.811ab3495dff9424ce59152fc652fd16cd023d27.de-de.xlf"
$ws.Range("K2").Value = "2016-08-31 05:40:55"

$ws.Range("J3").Value = "d2abd64f-9c0f-4af3-acf8-f4f5f601c048.05a88b359f114698fe9368b5815e5d2a3493305a.de-de.xlf"
$ws.Range("K3").Value = "2016-08-31 05:40:55"

# Recreate all four hyperlinks (A2, I2, A3, I3) in reading order so the
# relationship ids line up the way Excel would renumber them.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), $urlMd32, "", "", "32e4e091-9df4-414b-bc4c-8236a1306f0c.md")
$ws.Hyperlinks.Add($ws.Range("I2"), $urlMd32, "", "", "32e4e091-9df4-414b-bc4c-8236a1306f0c.md")
$ws.Hyperlinks.Add($ws.Range("A3"), $urlMdD2, "", "", "d2abd64f-9c0f-4af3-acf8-f4f5f601c048.md")
$ws.Hyperlinks.Add($ws.Range("I3"), $urlMdD2, "", "", "d2abd64f-9c0f-4af3-acf8-f4f5f601c048.md")

$ws.Range("I2").Font.Underline = $true
$ws.Range("I2").Font.Color = 15570276
$ws.Range("I3").Font.Underline = $true
$ws.Range("I3").Font.Color = 15570276

$ws.Columns.Item(3).ColumnWidth = 29.9777050018311
$ws.Columns.Item(9).ColumnWidth = 40
$ws.Columns.Item(10).ColumnWidth = 40

# ---------------- Overview ----------------
# Same status text swap (shared with the per-language sheets) plus the
# matching column widen (auto-fit side effect of the longer text).
$ws = $wb.Worksheets.Item("Overview")
$ws.Range("E2").Value = $statusText
$ws.Range("F2").Value = $statusText
$ws.Range("E3").Value = $statusText
$ws.Range("F3").Value = $statusText
$ws.Columns.Item(5).ColumnWidth = 29.9777050018311
$ws.Columns.Item(6).ColumnWidth = 29.9777050018311

Write-Host "Handback report generated."
